# Add a header row to the COMPASS_Synoptic_CB_MonMon_2025 sheet.
# 1. Insert a new blank row at the top, pushing all existing data down
#    by one row (this also naturally keeps the "Number" column's 1..115
#    sequence intact on rows 2..116, since the cell contents just move).
# 2. Populate the new row 1 with the column headers.
# 3. Move the active selection to F1 (matches the saved selection in the
#    target workbook).
# 4. Update the hidden _FilterDatabase defined name so it covers the new
#    data range (header excluded, same as before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Site"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Zone"
$ws.Range("E1").Value = "Replicate"
$ws.Range("F1").Value = "Depth"

$ws.Range("F1").Select()

$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=COMPASS_Synoptic_CB_MonMon_2025!`$A`$2:`$A`$126"
